{"js": "// Resume bullet edits:\n// 1) \"Designed and executed over 25,000 Email, Journey, and Automation test\n//     cases, ensuring seamless functionality of ... UTM parameters and\n//     fallbacks within each email, in accordance with the requestor's\n//     specifications.\"\n//    becomes\n//    \"Designed and executed thousands of manual test cases for Emails,\n//     Journeys, and Automations, ensuring seamless functionality of ...\n//     UTM parameters and fallbacks within each email.\"\n//\n// 2) \"Utilized SQL queries in Salesforce Marketing Cloud (SFMC) Query\n//     Studio ...\"\n//    becomes\n//    \"Utilized SQL queries and test scripts in Salesforce Marketing Cloud\n//     (SFMC) Query Studio ...\"\n\nconst body = context.document.body;\n\n// --- Edit 1a: \"over 25,000 Email, Journey, and Automation test cases, \" ---\nconst hit1 = body.search(\"over 25,000 Email, Journey, and Automation test cases, \", { matchCase: true });\nhit1.load(\"items\");\nawait context.sync();\nif (hit1.items.length > 0) {\n  hit1.items[0].insertText(\"thousands of manual test cases for Emails, Journeys, and Automations, \", \"Replace\");\n}\n\n// --- Edit 1b: drop the trailing \"in accordance with the requestor's specifications\" ---\nconst hit2 = body.search(\"within each email, in accordance with the requestor's specifications.\", { matchCase: true });\nhit2.load(\"items\");\nawait context.sync();\nif (hit2.items.length > 0) {\n  hit2.items[0].insertText(\"within each email.\", \"Replace\");\n}\n\n// --- Edit 2: \"Utilized SQL queries in Salesforce\" -> \"Utilized SQL queries and test scripts in Salesforce\" ---\nconst hit3 = body.search(\"Utilized SQL queries in Salesforce\", { matchCase: true });\nhit3.load(\"items\");\nawait context.sync();\nif (hit3.items.length > 0) {\n  hit3.items[0].insertText(\"Utilized SQL queries and test scripts in Salesforce\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Resume bullet edits:\n# 1) \"Designed and executed over 25,000 Email, Journey, and Automation test\n#     cases, ensuring seamless functionality of ... UTM parameters and\n#     fallbacks within each email, in accordance with the requestor's\n#     specifications.\"\n#    becomes\n#    \"Designed and executed thousands of manual test cases for Emails,\n#     Journeys, and Automations, ensuring seamless functionality of ...\n#     UTM parameters and fallbacks within each email.\"\n#\n# 2) \"Utilized SQL queries in Salesforce Marketing Cloud (SFMC) Query\n#     Studio ...\"\n#    becomes\n#    \"Utilized SQL queries and test scripts in Salesforce Marketing Cloud\n#     (SFMC) Query Studio ...\"\n\n$d = $word.ActiveDocument\n\n# --- Edit 1a: \"over 25,000 Email, Journey, and Automation test cases, \" ---\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Execute(\n    \"over 25,000 Email, Journey, and Automation test cases, \",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"thousands of manual test cases for Emails, Journeys, and Automations, \",\n    2\n)\n\n# --- Edit 1b: drop the trailing \"in accordance with the requestor's specifications\" ---\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\n    \"within each email, in accordance with the requestor's specifications.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"within each email.\",\n    2\n)\n\n# --- Edit 2: \"Utilized SQL queries in Salesforce\" -> \"Utilized SQL queries and test scripts in Salesforce\" ---\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Execute(\n    \"Utilized SQL queries in Salesforce\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Utilized SQL queries and test scripts in Salesforce\",\n    2\n)\n"}
